$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2's formatting down into row 3 (new "results" row) so the
# date/percentage number formats match without minting new style entries.
$ws.Range("A2:W2").Copy($ws.Range("A3:W3"))

# Now overwrite row 3 with the new ticker result values.
$ws.Range("A3").Value = 42632.883506944447
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 11353
$ws.Range("F3").Value = 615
$ws.Range("G3").Value = 64
$ws.Range("H3").Value = 35
$ws.Range("I3").Value = 84
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 15821
$ws.Range("L3").Value = 134
$ws.Range("M3").Value = 73
$ws.Range("N3").Value = 11
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = "Named"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.87
$ws.Range("S3").Value = 0.0351
$ws.Range("T3").Value = -2.08
$ws.Range("U3").Value = 15.16
$ws.Range("V3").Value = "N/A"
$ws.Range("W3").Value = 0
